$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.284.31'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.87%  '
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = '1.710.72'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.85%  '
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = '''1.003'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = '''224.90'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.39%  '
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = '''0.5291'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = '''1.003'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = '''0.06689'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +1.40%  '
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = '  -0.33%  '
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = '''20.86'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -3.63%  '
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = '''0.07690'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.31%  '
$ws.Range("E11").Style = "Normal"
$ws.Range("E12").Value = '  -2.45%  '
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = '1.945.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -0.86%  '
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = '1.703.60'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.28%  '
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = '''0.5856'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.52%  '
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = '0.0₅8219'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.81%  '
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = '''67.88'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.02%  '
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = '27.323.08'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.78%  '
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = '''221.76'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.31%  '
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = '  -0.01%  '
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = '''4.649'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -1.68%  '
$ws.Range("E21").Style = "Normal"
$ws.Range("E22").Value = '  -1.35%  '
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = '''6.029'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.99%  '
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = '''1.004'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.05%  '
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = '''144.88'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = '''1.693'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.14%  '
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = '  -2.28%  '
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = '''16.25'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.86%  '
$ws.Range("E29").Style = "Normal"
$ws.Range("E30").Value = '  -3.43%  '
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = '''1.292'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.97%  '
$ws.Range("E31").Style = "Normal"
$ws.Range("E32").Value = '  -2.80%  '
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = '''3.425'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.50%  '
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = '''1.629'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Value = '''2.870'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.46%  '
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = '''0.9560'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.99%  '
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = '''2.392'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.43%  '
$ws.Range("E37").Style = "Normal"
$ws.Range("D38").Value = '''0.5859'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -2.08%  '
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = '1.145.93'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.47%  '
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = '''0.01639'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = '''5.787'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.00%  '
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = '''1.004'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.02%  '
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = '''0.8384'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -1.53%  '
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = '''100.90'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.60%  '
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = '1.853.03'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -0.88%  '
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = '  -2.80%  '
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = '''57.58'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -2.32%  '
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = '''0.4568'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +2.21%  '
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = '''1.002'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.07%  '
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = '''8.103'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.83%  '
$ws.Range("E50").Style = "Normal"
$ws.Range("E51").Value = '  -0.79%  '
$ws.Range("E51").Style = "Normal"
